# Fill Hoja1 (sheet1) with the full list of year file entries (2002-2023),
# mirroring the data already present on Hoja2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2002.xlsx", "2002"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2003.xlsx", "2003"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2004.xlsx", "2004"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2005.xlsx", "2005"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2006.xlsx", "2006"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2007.xlsx", "2007"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2008.xlsx", "2008"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2009.xlsx", "2009"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2010.xlsx", "2010"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2011.xlsx", "2011"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2012.xlsx", "2012"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2013.xlsx", "2013"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2014.xlsx", "2014"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2015.xlsx", "2015"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2016.xlsx", "2016"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2017.xlsx", "2017"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2018.xlsx", "2018"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2019.xlsx", "2019"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2020.xlsx", "2020"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2021.xlsx", "2021"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2022.xlsx", "2022"),
    @("C:\Users\zaka\Desktop\MOTOGP\Excels\data\2023.xlsx", "2023")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("A2:B23").Select()
